# Insert a new row at position 726 (shifts rows 726:839 down to 727:840)
# and populate it with the new weekly data point, matching formatting of
# the surrounding rows (especially the date style on column D).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("726:726").Insert()

$ws.Range("A726").Value = 6
$ws.Range("B726").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C726").Value = "Metropolitana"
$ws.Range("D726").Value = 45218
$ws.Range("E726").Value = 13
$ws.Range("F726").Value = 100112012
$ws.Range("G726").Value = "Espinaca"
$ws.Range("H726").Value = "Sin especificar"
$ws.Range("I726").Value = "Primera"
$ws.Range("J726").Value = 510
$ws.Range("K726").Value = 7000
$ws.Range("L726").Value = 8000
$ws.Range("M726").Value = 7549
$ws.Range("N726").Value = '$/cuna 10 kilos'
$ws.Range("O726").Value = "Región Metropolitana"
$ws.Range("P726").Value = 755
$ws.Range("Q726").Value = 10
$ws.Range("R726").Value = "Hortaliza"
